# Generate Report for Handback
#
# The localization-status report is regenerated: both the zh-cn and de-de
# rows have been handed back and are now in sync with en-US, so:
#   - Status ("Ready for handoff") becomes "Handed back: in sync with en-US"
#     on the Overview sheet (columns E/F) and on the zh-cn / de-de sheets
#     (column C, "Status").
#   - The "Latest Handback DateTime" (column K) on the zh-cn / de-de sheets
#     is refreshed to the new handback timestamp.
#   - The "Error Detail" (column P) on the zh-cn / de-de sheets, which used
#     to warn that the handback file was stale, is cleared now that the
#     handback is current.
#   - Columns widened/narrowed to fit the new (longer) status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet row 2 ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-30 04:49:14"
$zhcn.Range("P2").Value = ""

# --- de-de sheet row 2 ---
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-30 04:49:20"
$dede.Range("P2").Value = ""

# --- Column width adjustments (status text is now longer) ---
# Target stored widths are 29.9777047293527 / 13.7470528738839; ColumnWidth
# gets quantized to the nearest 1/6 character unit by the host, so feed it
# the input that lands on the closest achievable stored width.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668   # -> stored 30
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668   # -> stored 30

$zhcn.Columns.Item(3).ColumnWidth  = 29.166666666666668      # -> stored 30
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334      # -> stored 13.6667

$dede.Columns.Item(3).ColumnWidth  = 29.166666666666668      # -> stored 30
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334      # -> stored 13.6667
